# Helper: convert an "RRGGBB" hex string into the packed Long value that the
# PowerPoint COM object model's RGBColor.RGB property expects (classic VBA
# RGB() encoding: value = R + G*256 + B*65536).
function ConvertTo-VbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Slide 5: re-style the table with the new table-style GUID ---------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{F98CD05B-CCEC-4A68-A7DC-2B4E70FDE13E}")

# --- 2. Swap the deck's theme palette from "Integral" (Red Violet) to the
#        stock "Office Theme" colors -------------------------------------
$officeColors = @(
    '000000',  # dk1
    'FFFFFF',  # lt1
    '44546A',  # dk2
    'E7E6E6',  # lt2
    '5B9BD5',  # accent1
    'ED7D31',  # accent2
    'A5A5A5',  # accent3
    'FFC000',  # accent4
    '4472C4',  # accent5
    '70AD47',  # accent6
    '0563C1',  # hlink
    '954F72'   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-VbaRgb($officeColors[$i - 1])
}
